# Add 5 "Dummy" data columns (H:L) with a header row and 43 rows of
# numeric data (rows 2-44), then update the sheet view (zoom + selection)
# to match the uploaded workbook revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells H1:L1
$headers = @("Dummy 1", "Dummy 2", "Dummy 3", "Dummy 4", "Dummy 5")
for ($j = 0; $j -lt $headers.Length; $j++) {
    $ws.Cells.Item(1, $j + 8).Value = $headers[$j]
}

# New data for H2:L44 (one row per existing data row)
$data = @(
    @(11,40,39,31,30),
    @(25,12,10,31,20),
    @(11,10,25,40,24),
    @(15,29,10,38,29),
    @(16,40,13,23,28),
    @(33,16,17,32,21),
    @(10,29,18,37,26),
    @(11,14,13,37,24),
    @(31,25,23,36,30),
    @(31,38,32,35,19),
    @(22,37,37,19,20),
    @(26,25,24,23,38),
    @(28,13,24,12,11),
    @(21,30,40,23,23),
    @(38,26,32,21,40),
    @(27,27,36,26,27),
    @(26,28,20,20,14),
    @(20,28,32,16,22),
    @(34,19,18,18,19),
    @(22,31,15,39,21),
    @(37,20,16,29,22),
    @(32,12,16,27,29),
    @(26,22,16,25,20),
    @(29,29,20,12,35),
    @(19,23,32,36,40),
    @(37,27,39,32,20),
    @(18,22,37,16,19),
    @(12,10,13,36,34),
    @(16,12,33,40,22),
    @(12,29,23,27,13),
    @(19,29,33,24,23),
    @(31,20,30,33,30),
    @(40,23,24,12,17),
    @(10,31,19,37,38),
    @(16,14,18,39,16),
    @(20,14,27,26,37),
    @(10,23,32,29,21),
    @(38,12,36,26,18),
    @(38,36,15,29,34),
    @(21,18,21,29,12),
    @(40,24,40,14,32),
    @(29,30,22,31,16),
    @(18,33,27,25,17)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 8).Value = $rowVals[$j]
    }
}

# Match the updated sheet view: zoom to 130% and select M5
$excel.ActiveWindow.Zoom = 130
$ws.Range("M5").Select()
